$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix capitalization of unit strings in column D (unit column)
# D14 corresponds to cap_Delta_17O -> unit should be "perMeg"
# D15 corresponds to d17O -> unit should be "perMil"
# D16 corresponds to d18O -> unit should be "perMil"
$ws.Range("D14").Value = "perMeg"
$ws.Range("D15").Value = "perMil"
$ws.Range("D16").Value = "perMil"

# Update the active selection to D16
$ws.Range("D16").Select()
